# Inventory: disable default normalization and fix parser robustness
#
# Beer Cost sheet: replace the single placeholder data row with the full
# (alphabetised) beer inventory count.
# Wine Cost sheet: populate the previously-empty sheet with a header row
# (styled to match the other inventory tabs) plus the wine inventory count.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Beer Cost
# ---------------------------------------------------------------------
$beer = $wb.Worksheets.Item("Beer Cost")

$beerItems = @(
    @("Athletic Upside Dawn Golden 12oz Can", 322),
    @("Athletic Upside Dawn NA (cans)", 41),
    @("Coors Light 12oz Can", 12),
    @("Fairhope I Drink Therefore I Amber Keg", 914.25),
    @("Grayton 30A Beach Blonde Keg", 63),
    @("Grayton Beach 30A Rosé 12oz Can", 1.5),
    @("High Noon (cans)", 20),
    @("High Rise Blood Orange 12oz Can", 27),
    @("High Rise Pineapple THC Seltzer (cans)", 13),
    @("Michelob Ultra (cans)", 68),
    @("Miller Lite (cans)", 40),
    @("Miller Lite 12oz Can", 205.5),
    @("Mom Water Linda 12oz Can", 57.5),
    @("Yuengling (cans)", 22)
)

$row = 2
foreach ($item in $beerItems) {
    $beer.Cells.Item($row, 1).Value = $item[0]
    $beer.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Wine Cost
# ---------------------------------------------------------------------
$wine = $wb.Worksheets.Item("Wine Cost")

$wine.Range("A1").Value = "Item"
$wine.Range("B1").Value = "Count"

# Match the bold/centered/bordered header style already used on the other
# inventory tabs (e.g. Beer Cost!A1:B1) instead of inventing a new style.
$beer.Range("A1:B1").Copy()
$wine.Range("A1:B1").PasteSpecial(-4122)

$wineItems = @(
    @("Carter's Lot Rosé of Pinot Noir 750ml", 62.5),
    @("The Shaker Red Blend", 1)
)

$row = 2
foreach ($item in $wineItems) {
    $wine.Cells.Item($row, 1).Value = $item[0]
    $wine.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
